$d = $word.ActiveDocument

$replacements = @(
    @("2025-01-07 Tuesday", "2025-01-08 Wednesday"),
    @("124×6=744", "802×6=4812"),
    @("996×3=2988", "103×2=206"),
    @("321×3=963", "330×9=2970"),
    @("132×7=924", "679×8=5432"),
    @("571×8=4568", "566×5=2830"),
    @("357×9=3213", "736×6=4416"),
    @("386×4=1544", "622×5=3110"),
    @("759×6=4554", "560×4=2240"),
    @("624×6=3744", "730×5=3650"),
    @("744×2=1488", "716×9=6444"),
    @("409×3=1227", "140×8=1120"),
    @("918×6=5508", "827×2=1654"),
    @("675×9=6075", "641×4=2564"),
    @("356×7=2492", "463×7=3241"),
    @("739×6=4434", "981×6=5886"),
    @("631×4=2524", "900×6=5400"),
    @("452×9=4068", "407×9=3663"),
    @("402×5=2010", "955×6=5730"),
    @("681×8=5448", "581×6=3486"),
    @("417×3=1251", "359×4=1436"),
    @("823×6=4938", "824×2=1648"),
    @("317×4=1268", "931×5=4655"),
    @("242×3=726", "298×5=1490"),
    @("861×9=7749", "214×4=856"),
    @("596×8=4768", "892×9=8028")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
